# Apply the edit described by the diff:
#  - Sheet10 (the one with dimension A1:S3, selection J2:J5) -> selection becomes A2
#  - Sheet21 -> selection becomes whole row 1 (A1:XFD1), no longer the active tab
#  - New Sheet22 appended at the end, with data + becomes the active tab,
#    selection B5
#  - The new text written to Sheet22!R4 is entered with a leading apostrophe
#    (quote-prefix) exactly like a user typing '-100, 0, 100 into Excel, which
#    is what produces the new quotePrefix cellXfs entry.

$wb = $excel.ActiveWorkbook

# --- Sheet10: A1:S3 / old selection J2:J5 -> new selection A2 ---------------
$ws10 = $wb.Worksheets.Item(10)
[void]$ws10.Range("A2").Select()

# --- Sheet21: old selection J3:J5 (tabSelected) -> whole row 1 selected, ----
# --- and it stops being the active tab (Sheet22 takes over below) ----------
$ws21 = $wb.Worksheets.Item(21)
[void]$ws21.Range("A1:XFD1").Select()

# --- Add the new Sheet22 right after the current last sheet (Sheet21) ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws22 = $wb.Worksheets.Add($null, $lastSheet)
$ws22.Name = "Sheet22"

# Header row (row 1) -- same 19 headers used on every other sheet
$ws22.Range("A1").Value = "serija"
$ws22.Range("B1").Value = "enota"
$ws22.Range("C1").Value = "legenda"
$ws22.Range("D1").Value = "barva"
$ws22.Range("E1").Value = "tip"
$ws22.Range("F1").Value = "stacked"
$ws22.Range("G1").Value = "drseca_obdobja"
$ws22.Range("H1").Value = "drseca_poravnava"
$ws22.Range("I1").Value = "rast"
$ws22.Range("J1").Value = "indeks_obdobje"
$ws22.Range("K1").Value = "velikost"
$ws22.Range("L1").Value = "naslov"
$ws22.Range("M1").Value = "xmin"
$ws22.Range("N1").Value = "xmax"
$ws22.Range("O1").Value = "opomba"
$ws22.Range("P1").Value = "stolpci_legende"
$ws22.Range("Q1").Value = "datum_podatkov"
$ws22.Range("R1").Value = "leva_y_os"
$ws22.Range("S1").Value = "desna_y_os"

# Row 2 -- written A, C, B (matches the order new shared strings were added)
$ws22.Range("A2").Value = "MF--OB--001--7--A"
$ws22.Range("C2").Value = "SKUPAJ PRIHODKI -- Letno"
$ws22.Range("B2").Value = "eur"

# Row 3 -- written A, C, B
$ws22.Range("A3").Value = "MF--OB--001--7--M"
$ws22.Range("C3").Value = "SKUPAJ PRIHODKI -- Mesečno"
$ws22.Range("B3").Value = "eur"

# Row 4 -- written A, C, J, R, S
$ws22.Range("A4").Value = "SURS--0300230S--P31_S14_D--G4--N--Q"
$ws22.Range("C4").Value = "19 Proizvodnja koksa in naftnih derivatov -- Podjetja v % -- negotove gospodarske razmere"
$ws22.Range("J4").Value = 2015
# Leading apostrophe -> Excel quote-prefix (text-forced) entry, producing the
# new cellXfs style (quotePrefix="1") applied only to this one cell.
$ws22.Range("R4").Value = "'-100, 0, 100"
$ws22.Range("S4").Value = "0, 50000, 100000"

# Make Sheet22 the active tab with the matching selection, as the very last
# selection action in the script (so it "wins" the tabSelected flag).
[void]$ws22.Range("B5").Select()
